$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values to write. Values that look like plain numbers must be forced
# to remain text (matching the source inline-string data), so for those we
# temporarily switch the cell to Text format, assign, then restore the style.
$textForcedValues = @{
    'G2' = '1.15'
    'F3' = '98'
    'G3' = '1.18'
    'G4' = '1.70'
    'F5' = '95'
    'G5' = '1.10'
    'F6' = '89'
    'G6' = '1.40'
    'F7' = '94'
    'G7' = '1.44'
    'G8' = '1.83'
    'F9' = '96'
    'G9' = '1.00'
    'F10' = '84'
    'G10' = '2.25'
    'F11' = '56'
    'G11' = '3.60'
    'F12' = '93'
    'G12' = '1.25'
    'F13' = '71'
    'G13' = '3.05'
    'F14' = '77'
    'G14' = '3.15'
    'F15' = '50'
    'G15' = '4.75'
    'F16' = '100'
    'G16' = '1.25'
    'G17' = '1.65'
    'G18' = '1.33'
    'F19' = '47'
    'G19' = '3.20'
    'F20' = '70'
    'G20' = '2.30'
    'F21' = '64'
    'G21' = '1.95'
    'G22' = '7.50'
    'G23' = '5.00'
}

$plainValues = @{
    'G1' = 'Odds'
    'E2' = '52/54 Win Tips'
    'E3' = '46/47 Win Tips'
    'E4' = '43/48 Win Tips'
    'E5' = '42/44 Win Tips'
    'E6' = '34/38 Win Tips'
    'E7' = '30/32 Win Tips'
    'E8' = '25/26 Win Tips'
    'E9' = '22/23 Win Tips'
    'E10' = '16/19 Win Tips'
    'A11' = 'Sporting Gijon v Racing Santander'
    'B11' = 'Racing Santander'
    'C11' = 'Spain Segunda'
    'D11' = '2025-10-12T14:15:00.000Z'
    'E11' = '14/25 Win Tips'
    'A12' = 'Egypt v Guinea Bissau'
    'B12' = 'Egypt'
    'C12' = 'W Cup African Qual'
    'D12' = '2025-10-12T19:00:00.000Z'
    'E12' = '13/14 Win Tips'
    'A13' = 'Burgos v Valladolid'
    'B13' = 'Draw'
    'D13' = '2025-10-12T16:30:00.000Z'
    'E13' = '12/17 Win Tips'
    'A14' = 'Chad v Central African Republic'
    'C14' = 'W Cup African Qual'
    'D14' = '2025-10-12T16:00:00.000Z'
    'E14' = '10/13 Win Tips'
    'A15' = 'Zambia v Niger'
    'B15' = 'Niger'
    'D15' = '2025-10-12T13:00:00.000Z'
    'E15' = '9/18 Win Tips'
    'A16' = 'Burkina Faso v Ethiopia'
    'B16' = 'Burkina Faso'
    'D16' = '2025-10-12T19:00:00.000Z'
    'E16' = '9/9 Win Tips'
    'A17' = 'Malta v Bosnia Herzegovina'
    'B17' = 'Bosnia Herzegovina'
    'C17' = 'International'
    'D17' = '2025-10-12T17:00:00.000Z'
    'A18' = 'Ghana v Comoros'
    'B18' = 'Ghana'
    'E19' = '7/15 Win Tips'
    'A20' = 'Malaga v Deportivo La Coruna'
    'B20' = 'Deportivo La Coruna'
    'C20' = 'Spain Segunda'
    'D20' = '2025-10-12T19:00:00.000Z'
    'E20' = '7/10 Win Tips'
    'A21' = 'Nottm Forest v Chelsea'
    'B21' = 'Chelsea'
    'C21' = 'England Premier League'
    'D21' = '2025-10-18T11:30:00.000Z'
    'E21' = '7/11 Win Tips'
}

foreach ($addr in $plainValues.Keys) {
    $ws.Range($addr).Value = $plainValues[$addr]
}

foreach ($addr in $textForcedValues.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedValues[$addr]
    $cell.Style = $origStyle
}
